$d = $word.ActiveDocument

# 1) Merge the two runs in the "Aqui já será outra ideia ... " paragraph into
#    a single run carrying the full, combined text (same run formatting).
#    The two original runs share identical rPr, so replacing the whole
#    (cross-run) match with the same text collapses them into one run.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Aqui já será outra ideia … Salvando os trechos de códigos desenvolvidos",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Aqui já será outra ideia … Salvando os trechos de códigos desenvolvidos",
    2) | Out-Null

# Locate the paragraph that now holds that merged sentence, so the two
# following (previously unformatted) paragraphs can be found relative to it.
$target = $find.Parent.Paragraphs.First.Index

# 2) The empty paragraph right after it gets the same run formatting
#    (Segoe UI, non-bold, italic, size 12/24 half-points).
$p1 = $d.Paragraphs($target + 1)
$r1 = $p1.Range
$r1.Font.Name = "Segoe UI"
$r1.Font.NameAscii = "Segoe UI"
$r1.Font.Bold = $false
$r1.Font.BoldBi = $false
$r1.Font.Italic = $true
$r1.Font.ItalicBi = $true
$r1.Font.Size = 12
$r1.Font.SizeBi = 12

# 3) The following paragraph (containing the anchored picture) gets the same
#    run formatting applied to its run too.
$p2 = $d.Paragraphs($target + 2)
$r2 = $p2.Range
$r2.Font.Name = "Segoe UI"
$r2.Font.NameAscii = "Segoe UI"
$r2.Font.Bold = $false
$r2.Font.BoldBi = $false
$r2.Font.Italic = $true
$r2.Font.ItalicBi = $true
$r2.Font.Size = 12
$r2.Font.SizeBi = 12

Write-Output ("target=" + $target)
